$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period sequence (1705 -> 1803) for "Periodo Mora" column (E16:E26)
$periodos = @("1705", "1706", "1707", "1708", "1709", "1710", "1711", "1712", "1801", "1802", "1803")

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("G$row").Value = 781242
}
